$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Wins, Losses, Ties in AC1:AE1, matching the bold/centered/bordered
# style already used by the rest of the row-1 header cells (copy format from AB1).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Season record repeated for every player row (2-38): Wins=84, Losses=78, Ties=0
$ws.Range("AC2:AC38").Value = 84
$ws.Range("AD2:AD38").Value = 78
$ws.Range("AE2:AE38").Value = 0
